$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$yellow = 65535

# ---------------------------------------------------------------------------
# 1) Build up the cellXfs / fonts / fills table in the exact order the
#    target workbook uses:
#      idx1 = fill (yellow) only
#      idx2 = bold font only
#      idx3 = left-aligned only
#      idx4 = bold font + fill (yellow)
# ---------------------------------------------------------------------------

# NOTE: this host's Range(..) does not reliably apply formatting to every
# area of a comma-separated union (only the first area "sticks" once a
# .Value is written afterwards) -- so every range below is applied
# individually rather than as a union.

# idx1: fill-only -> the "blank filler" cells to the right of each section title
$ws.Range("B1:R1").Interior.Color = $yellow
$ws.Range("B8:P8").Interior.Color = $yellow
$ws.Range("B15:P15").Interior.Color = $yellow

# idx2: bold-only -> the 3 column-header rows
$ws.Range("A2:E2").Font.Bold = $true
$ws.Range("A9:E9").Font.Bold = $true
$ws.Range("A16:E16").Font.Bold = $true

# idx3: left-align only -> the numbered serial cells
$ws.Range("A3").HorizontalAlignment = -4131
$ws.Range("A10").HorizontalAlignment = -4131
$ws.Range("A17").HorizontalAlignment = -4131

# idx4: bold + fill -> the 3 section title cells
$ws.Range("A1").Font.Bold = $true
$ws.Range("A1").Interior.Color = $yellow
$ws.Range("A8").Font.Bold = $true
$ws.Range("A8").Interior.Color = $yellow
$ws.Range("A15").Font.Bold = $true
$ws.Range("A15").Interior.Color = $yellow

# ---------------------------------------------------------------------------
# 2) Cell content. Row 1-6 existing block: re-point A3/B3 values.
# ---------------------------------------------------------------------------
$ws.Range("A1").Value = "Grocery"

$ws.Range("A3").Value = 1
$ws.Range("B3").Value = "Grocery"

# ---------------------------------------------------------------------------
# 3) New "Mobile" block (rows 8-13). Values are written in the precise
#    order needed so the shared-strings table lands in the same order as
#    the target workbook (column by column, with a couple of cells the
#    original author circled back and filled in last).
# ---------------------------------------------------------------------------
$ws.Range("A8").Value = "Mobile"

$ws.Range("A9").Value = "TestScenarioid"
$ws.Range("B9").Value = "TestScenario"
$ws.Range("C9").Value = "Rtm id"
$ws.Range("D9").Value = "TestCaseid"
$ws.Range("E9").Value = "TestCase"

$ws.Range("A10").Value = 2
$ws.Range("B10").Value = "Mobile"

$ws.Range("C10").Value = "R5"
$ws.Range("C11").Value = "R6"
$ws.Range("C12").Value = "R7"
$ws.Range("C13").Value = "R8"

$ws.Range("D10").Value = "TC5"
$ws.Range("D11").Value = "TC6"
$ws.Range("D12").Value = "TC7"
$ws.Range("D13").Value = "TC8"

$ws.Range("E10").Value = "there should be ratings of the mobile"
$ws.Range("E12").Value = "features option "
$ws.Range("E13").Value = "camera clarity should be mentioned"

# ---------------------------------------------------------------------------
# 4) New "Appliance Selection" block (rows 15-20).
# ---------------------------------------------------------------------------
$ws.Range("A15").Value = "Appliance Selection"

$ws.Range("A16").Value = "TestScenarioid"
$ws.Range("B16").Value = "TestScenario"
$ws.Range("C16").Value = "Rtm id"
$ws.Range("D16").Value = "TestCaseid"
$ws.Range("E16").Value = "TestCase"

$ws.Range("A17").Value = 3
$ws.Range("B17").Value = "Appliance"

$ws.Range("C17").Value = "R9"
$ws.Range("C18").Value = "R10"
$ws.Range("C19").Value = "R11"
$ws.Range("C20").Value = "R12"

$ws.Range("D17").Value = "TC9"
$ws.Range("D18").Value = "TC10"
$ws.Range("D19").Value = "TC11"
$ws.Range("D20").Value = "TC12"

$ws.Range("E18").Value = "when clicked arrow mark of right that should be navigated "
$ws.Range("E19").Value = "back to top option"
$ws.Range("E20").Value = "when clicked price hightolow that should be worked"

# These two were edited last by the original author (their shared-string
# ids land at the very end of the table), so they're set last here too.
$ws.Range("E17").Value = "when scrolled down that should be scrolled"
$ws.Range("E11").Value = "when price is clicked someother page should be opened"

# ---------------------------------------------------------------------------
# 5) Column E width + final selection.
# ---------------------------------------------------------------------------
$ws.Columns.Item(5).ColumnWidth = 49.15

$ws.Range("A15").Select()
